$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Move the "_GoBack" bookmark from the (now empty) paragraph that
#    used to sit right before "Payer->UA: Select Payment Instrument"
#    to the very start of the document, immediately before the
#    "@startuml" run. Re-adding a bookmark with the same name moves
#    it (Word enforces unique bookmark names), so the old location is
#    cleared automatically. The paragraph that is left empty by the
#    move is then deleted so the two adjacent empty paragraphs merge
#    back into one, exactly as in the target document.
# ------------------------------------------------------------------

# Inserting the bookmark directly at offset 0 has a boundary quirk in
# this host, so nudge it in: insert a throwaway character, anchor the
# (zero-length) bookmark range right after it, then remove the
# throwaway character again. The bookmark stays put at position 0.
$startRange = $d.Range(0, 0)
$startRange.InsertBefore("X")

$bmRange = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

$dummyRange = $d.Range(0, 1)
$dummyRange.Delete()

# Find the now-vacated bookmark paragraph (the empty paragraph right
# before "Payer->UA: Select Payment Instrument") and remove it.
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -eq "`r" -or $para.Range.Text -eq "") {
        $nextPara = $para.Next()
        if ($nextPara -ne $null -and $nextPara.Range.Text.StartsWith("Payer->UA: Select Payment Instrument")) {
            $para.Range.Delete()
            $found = $true
            break
        }
    }
}

# ------------------------------------------------------------------
# 2. Make the raw PlantUML source text hidden (the commit adds the
#    raw .pml text but keeps it invisible, relying on the rendered
#    diagram image instead).
# ------------------------------------------------------------------
$d.Styles("PlantUML").Font.Hidden = $true
$d.Styles("PlantUMLChar").Font.Hidden = $true
